# Convert the "answer" column (G) from letter values (A/B/C/D) to the
# corresponding numeric option index (1/2/3/4) for every question row,
# then leave the active selection on F2 (matching the author's edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map the textual answer letters to their numeric option index.
$letterToNumber = @{
    "A" = 1
    "B" = 2
    "C" = 3
    "D" = 4
}

# Find the last used row in column A (question index column) so this
# works regardless of exact row count; data starts on row 2 (row 1 is
# the header: id/text/option1/option2/option3/option4/answer).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $letter = [string]$cell.Value2
    if ($letterToNumber.ContainsKey($letter)) {
        $cell.Value2 = $letterToNumber[$letter]
    }
}

# Match the author's final selection (F2) after the edit.
$ws.Range("F2").Select()
